# Add season-record columns (Wins / Losses / Ties) to the DET_1994 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the formatting of an existing header cell
# (bold, bordered, centered) onto the three new header cells so the style
# index matches the rest of the header row, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-38): the 1994 Detroit Tigers finished the (strike
# shortened) season 53-62-0, so every player row gets the same team record.
$wins = 53
$losses = 62
$ties = 0

for ($row = 2; $row -le 38; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
